# Update cryptocurrency price and 1h volume-change figures to the latest
# scraped values (GitHub Actions refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $text) {
    # Force the cell to be stored as text (matching the original inline-string
    # cells) even when the new value looks like a plain number, then drop the
    # temporary text format so the cell keeps its original (default) style.
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.ClearFormats()
}

Set-TextValue $ws.Range("D2") "41.623.26"
Set-TextValue $ws.Range("E2") "  -1.02%  "

Set-TextValue $ws.Range("D3") "2.226.69"
Set-TextValue $ws.Range("E3") "  -0.99%  "

Set-TextValue $ws.Range("E4") "  -0.09%  "

Set-TextValue $ws.Range("D5") "252.00"
Set-TextValue $ws.Range("E5") "  +8.06%  "

Set-TextValue $ws.Range("D6") "0.627"
Set-TextValue $ws.Range("E6") "  -1.43%  "

Set-TextValue $ws.Range("D7") "71.01"
Set-TextValue $ws.Range("E7") "  +1.24%  "

Set-TextValue $ws.Range("E8") "  -0.04%  "

Set-TextValue $ws.Range("D9") "0.568"
Set-TextValue $ws.Range("E9") "  +1.40%  "

Set-TextValue $ws.Range("D10") "42.34"
Set-TextValue $ws.Range("E10") "  +18.00%  "

Set-TextValue $ws.Range("E11") "  -3.16%  "

Set-TextValue $ws.Range("D12") "58.77"
Set-TextValue $ws.Range("E12") "  +0.56%  "

Set-TextValue $ws.Range("E13") "  +0.29%  "

Set-TextValue $ws.Range("D14") "7.02"
Set-TextValue $ws.Range("E14") "  +2.75%  "

Set-TextValue $ws.Range("D15") "2.556.54"
Set-TextValue $ws.Range("E15") "  -0.90%  "

Set-TextValue $ws.Range("D16") "14.90"
Set-TextValue $ws.Range("E16") "  -1.82%  "

Set-TextValue $ws.Range("D17") "0.853"
Set-TextValue $ws.Range("E17") "  -1.45%  "

Set-TextValue $ws.Range("D18") "2.225.68"
Set-TextValue $ws.Range("E18") "  -0.77%  "

Set-TextValue $ws.Range("D19") "41.555.33"
Set-TextValue $ws.Range("E19") "  -0.95%  "

Set-TextValue $ws.Range("D20") "0.0₃0967"
Set-TextValue $ws.Range("E20") "  -1.67%  "

Set-TextValue $ws.Range("E21") "  -1.16%  "

Set-TextValue $ws.Range("D22") "72.84"
Set-TextValue $ws.Range("E22") "  -1.08%  "

Set-TextValue $ws.Range("D23") "2.27"
Set-TextValue $ws.Range("E23") "  +11.00%  "

Set-TextValue $ws.Range("D24") "234.60"
Set-TextValue $ws.Range("E24") "  -1.46%  "

Set-TextValue $ws.Range("D25") "3.86"
Set-TextValue $ws.Range("E25") "  +6.02%  "

Set-TextValue $ws.Range("D26") "0.999"
Set-TextValue $ws.Range("E26") "  -0.01%  "

Set-TextValue $ws.Range("E27") "  +5.80%  "

Set-TextValue $ws.Range("D28") "10.48"
Set-TextValue $ws.Range("E28") "  +4.02%  "

Set-TextValue $ws.Range("E29") "  +1.27%  "

Set-TextValue $ws.Range("D30") "171.41"
Set-TextValue $ws.Range("E30") "  +1.06%  "

Set-TextValue $ws.Range("D31") "20.60"
Set-TextValue $ws.Range("E31") "  -0.45%  "

Set-TextValue $ws.Range("E32") "  +0.89%  "

Set-TextValue $ws.Range("E33") "  -2.12%  "

Set-TextValue $ws.Range("D34") "5.57"
Set-TextValue $ws.Range("E34") "  +1.00%  "

Set-TextValue $ws.Range("E35") "  -0.01%  "

Set-TextValue $ws.Range("D36") "26.65"
Set-TextValue $ws.Range("E36") "  +19.97%  "

Set-TextValue $ws.Range("E37") "  -2.50%  "

Set-TextValue $ws.Range("D38") "3.99"
Set-TextValue $ws.Range("E38") "  +10.19%  "

Set-TextValue $ws.Range("E39") "  +9.56%  "

Set-TextValue $ws.Range("E40") "  +0.72%  "

Set-TextValue $ws.Range("D41") "69.54"
Set-TextValue $ws.Range("E41") "  +3.25%  "

Set-TextValue $ws.Range("D42") "5.99"
Set-TextValue $ws.Range("E42") "  -1.19%  "

Set-TextValue $ws.Range("D43") "12.37"
Set-TextValue $ws.Range("E43") "  +23.75%  "

Set-TextValue $ws.Range("E44") "  +9.27%  "

Set-TextValue $ws.Range("D45") "5.06"
Set-TextValue $ws.Range("E45") "  +2.62%  "

Set-TextValue $ws.Range("E46") "  +9.89%  "

Set-TextValue $ws.Range("E47") "  -3.33%  "

Set-TextValue $ws.Range("E48") "  +0.90%  "

Set-TextValue $ws.Range("E49") "  -0.02%  "

Set-TextValue $ws.Range("E50") "  +6.15%  "

Set-TextValue $ws.Range("E51") "  +1.65%  "
